$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BigWigs")
$ws.Activate()

# Rows 74-77: Status column (D) flips from "needs testing" to "Working - test"
$ws.Range("D74:D77").Value = "Working - test"

# New rows 99-100: Molten Core / Lucifron mechanics
$ws.Cells.Item(99, 1).Value = "Molten Core"
$ws.Cells.Item(99, 2).Value = "Lucifron"
$ws.Cells.Item(99, 3).Value = "Lucifron's Curse Timer"
$ws.Cells.Item(99, 4).Value = "Working - test"

$ws.Cells.Item(100, 1).Value = "Molten Core"
$ws.Cells.Item(100, 2).Value = "Lucifron"
$ws.Cells.Item(100, 3).Value = "Impending Doom Timer"
$ws.Cells.Item(100, 4).Value = "Working - test"

# Extend the Status data validation list down to the new last row
$ws.Range("D2:D100").Validation.Delete()
$ws.Range("D2:D100").Validation.Add(3, 1, 1, "=Status")

# Move the on-screen selection to match the saved view state
$ws.Range("D73").Select()
